$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the Index column (A2:A11) down by 1, so the sequence
# starts at 1 instead of 2 (rows r=2..11 hold A values 2..11 -> 1..10)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
